$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.742.85'
$ws.Range("E2").Value = '  +1.20%  '

$ws.Range("D3").Value = '3.283.20'
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.19%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  -0.59%  '

$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("E10").Value = '  -1.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.420'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.56%  '

$ws.Range("D12").Value = '3.851.81'
$ws.Range("E12").Value = '  +0.20%  '

$ws.Range("E13").Value = '  -0.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.56%  '

$ws.Range("D15").Value = '68.711.24'
$ws.Range("E15").Value = '  +1.17%  '

$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("D17").Value = '3.266.03'
$ws.Range("E17").Value = '  -0.40%  '

$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '395.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.71%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.519'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.00%  '

$ws.Range("E25").Value = '  +0.70%  '

$ws.Range("E26").Value = '  +4.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.68%  '

$ws.Range("E28").Value = '  -0.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.16%  '

$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.76%  '

$ws.Range("E32").Value = '  +3.21%  '

$ws.Range("E33").Value = '  +1.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("E35").Value = '  -0.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.66%  '

$ws.Range("E37").Value = '  +8.84%  '

$ws.Range("E38").Value = '  -3.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.59%  '

$ws.Range("E42").Value = '  -3.22%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0692'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.32%  '

$ws.Range("D46").Value = '2.660.96'
$ws.Range("E46").Value = '  -0.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '341.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.66%  '

$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.19%  '

$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '31.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.994'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.76%  '
